$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 381 (rows 381..428 shift down to 382..429,
# carrying their existing data/formatting with them).
$ws.Rows("381:381").Insert()

# Populate the newly inserted row 381 with the new observation.
$ws.Cells.Item(381, 1).Value = 6
$ws.Cells.Item(381, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(381, 3).Value = "Metropolitana"
$ws.Cells.Item(381, 4).Value = 45142
$ws.Cells.Item(381, 5).Value = 13
$ws.Cells.Item(381, 6).Value = 100112026
$ws.Cells.Item(381, 7).Value = "Haba"
$ws.Cells.Item(381, 8).Value = "Sin especificar"
$ws.Cells.Item(381, 9).Value = "Primera"
$ws.Cells.Item(381, 10).Value = 410
$ws.Cells.Item(381, 11).Value = 14000
$ws.Cells.Item(381, 12).Value = 15000
$ws.Cells.Item(381, 13).Value = 14439
$ws.Cells.Item(381, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(381, 15).Value = "Provincia de Copiapó"
$ws.Cells.Item(381, 16).Value = 578
$ws.Cells.Item(381, 17).Value = 25
$ws.Cells.Item(381, 18).Value = "Hortaliza"
